$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# The review from danfogel100@gmail.com / avishaybar12@gmail.com
# (row 2) is being removed. Deleting the entire row shifts every
# subsequent review up by one row, and shrinks the used range from
# A1:F20 down to A1:F19.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).Delete()

# ---------------------------------------------------------------------
# Row deletion does not renumber the worksheet's Hyperlinks collection
# in this environment, so rebuild it explicitly: drop every hyperlink,
# then re-add one per remaining "email"/"recovery" cell, pointing at
# the same mailto: addresses, now anchored to their shifted rows.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @(
    @{ Cell = "C2";  Email = "snirelad61@gmail.com" },
    @{ Cell = "C3";  Email = "eligitel@gmail.com" },
    @{ Cell = "D3";  Email = "ronenchen27@gmail.com" },
    @{ Cell = "C4";  Email = "zaittomer@gmail.com" },
    @{ Cell = "C5";  Email = "rontiddler560@gmail.com" },
    @{ Cell = "C6";  Email = "gregneri12@gmail.com" },
    @{ Cell = "C7";  Email = "snizzvered@gmail.com" },
    @{ Cell = "C8";  Email = "budoyoni2@gmail.com" },
    @{ Cell = "C10"; Email = "hermanliran@gmail.com" },
    @{ Cell = "C11"; Email = "gazittalia1@gmail.com" },
    @{ Cell = "D11"; Email = "hermanliran@gmail.com" },
    @{ Cell = "C12"; Email = "leviadlevi22@gmail.com" },
    @{ Cell = "D12"; Email = "gazittalia1@gmail.com" },
    @{ Cell = "C13"; Email = "freelancernachus@gmail.com" },
    @{ Cell = "C14"; Email = "nevilgreen@gmail.com" },
    @{ Cell = "D14"; Email = "vikicrestina@gmail.com" },
    @{ Cell = "C15"; Email = "veredsnir12@gmail.com" },
    @{ Cell = "D15"; Email = "kevinkors122@gmail.com" },
    @{ Cell = "C16"; Email = "stevewonder3001@gmail.com" },
    @{ Cell = "D16"; Email = "budoyoni@gmail.com" },
    @{ Cell = "C17"; Email = "stclerari834@gmail.com" },
    @{ Cell = "C18"; Email = "stcydouel274@gmail.com" },
    @{ Cell = "C19"; Email = "kevinkors122@gmail.com" },
    @{ Cell = "D19"; Email = "sinuspai@gmail.com" }
)

foreach ($link in $links) {
    $target = $ws.Range($link.Cell)
    $ws.Hyperlinks.Add($target, "mailto:" + $link.Email, "", "", $link.Email) | Out-Null
}

# ---------------------------------------------------------------------
# Restore the cursor position that was saved with the workbook.
# ---------------------------------------------------------------------
$ws.Range("E2").Select()
